$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: widen columns E and F (Latest Handoff / generate date area) ---
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet: widen columns C, I, J ---
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# Status column now reflects a completed handback
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

# Populate "Latest Target File" (I) and "Latest Handback File" (J) for row 2/3
$zhcn.Range("I2").Value = "0fea17c1-a2ec-43a1-9a95-dbaf8bd228db.md"
$zhcn.Range("J2").Value = "0fea17c1-a2ec-43a1-9a95-dbaf8bd228db.cb75e3916162626dff2518a3c442d16a8ede5e36.zh-cn.xlf"
$zhcn.Range("I3").Value = "2629fe80-0dd1-4b1f-8be5-e244c09476e4.md"
$zhcn.Range("J3").Value = "2629fe80-0dd1-4b1f-8be5-e244c09476e4.b0e41a56b80e7bdeffe5513400e120746b45521b.zh-cn.xlf"

# Latest Handback DateTime for zh-cn
$zhcn.Range("K2").Value = "2016-09-06 00:53:59"
$zhcn.Range("K3").Value = "2016-09-06 00:53:59"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dda7961b612a3563029f51fe72394ea9e1fc58e7/e2e/0fea17c1-a2ec-43a1-9a95-dbaf8bd228db.md", "", "", "0fea17c1-a2ec-43a1-9a95-dbaf8bd228db.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dda7961b612a3563029f51fe72394ea9e1fc58e7/e2e/2629fe80-0dd1-4b1f-8be5-e244c09476e4.md", "", "", "2629fe80-0dd1-4b1f-8be5-e244c09476e4.md")

# --- de-de sheet: widen columns C, I, J ---
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("I2").Value = "0fea17c1-a2ec-43a1-9a95-dbaf8bd228db.md"
$dede.Range("J2").Value = "0fea17c1-a2ec-43a1-9a95-dbaf8bd228db.cb75e3916162626dff2518a3c442d16a8ede5e36.de-de.xlf"
$dede.Range("I3").Value = "2629fe80-0dd1-4b1f-8be5-e244c09476e4.md"
$dede.Range("J3").Value = "2629fe80-0dd1-4b1f-8be5-e244c09476e4.b0e41a56b80e7bdeffe5513400e120746b45521b.de-de.xlf"

# Latest Handback DateTime for de-de (distinct timestamp from zh-cn)
$dede.Range("K2").Value = "2016-09-06 00:54:12"
$dede.Range("K3").Value = "2016-09-06 00:54:12"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dda7961b612a3563029f51fe72394ea9e1fc58e7/e2e/0fea17c1-a2ec-43a1-9a95-dbaf8bd228db.md", "", "", "0fea17c1-a2ec-43a1-9a95-dbaf8bd228db.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dda7961b612a3563029f51fe72394ea9e1fc58e7/e2e/2629fe80-0dd1-4b1f-8be5-e244c09476e4.md", "", "", "2629fe80-0dd1-4b1f-8be5-e244c09476e4.md")
